$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '54.678.85'
$ws.Range('E2').Value = '  +0.53%  '
$ws.Range('D3').Value = '2.287.03'
$ws.Range('E3').Value = '  +0.08%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '506.28'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.55%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '129.14'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.64%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.996'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('E8').Value = '  +0.11%  '
$ws.Range('D9').Value = '2.309.24'
$ws.Range('E9').Value = '  +0.56%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0970'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +1.44%  '
$ws.Range('E11').Value = '  +1.69%  '
$ws.Range('E12').Value = '  +2.55%  '
$ws.Range('E13').Value = '  +4.28%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '23.57'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +1.99%  '
$ws.Range('D15').Value = '2.696.87'
$ws.Range('E15').Value = '  +0.18%  '
$ws.Range('D16').Value = '54.692.00'
$ws.Range('E16').Value = '  +0.64%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.0000131'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +1.22%  '
$ws.Range('D18').Value = '2.256.95'
$ws.Range('E18').Value = '  -1.62%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '10.64'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +3.43%  '
$ws.Range('E20').Value = '  +1.30%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.67'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +3.84%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '308.11'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.73%  '
$ws.Range('E23').Value = '  -0.22%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '60.39'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -2.80%  '
$ws.Range('E25').Value = '  -0.57%  '
$ws.Range('E26').Value = '  -0.32%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.50'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +2.00%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '171.59'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -1.75%  '
$ws.Range('E29').Value = '  +1.68%  '
$ws.Range('E30').Value = '  +1.63%  '
$ws.Range('E31').Value = '  +0.23%  '
$ws.Range('E32').Value = '  +5.40%  '
$ws.Range('E34').Value = '  +1.10%  '
$ws.Range('E35').Value = '  -0.28%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.908'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -3.73%  '
$ws.Range('E37').Value = '  +0.56%  '
$ws.Range('E38').Value = '  +1.48%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '36.67'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +1.67%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.377'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.54%  '
$ws.Range('E41').Value = '  +1.13%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '133.01'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +6.23%  '
$ws.Range('E43').Value = '  +0.50%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '4.87'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.55%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '252.91'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +4.76%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0502'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +1.06%  '
$ws.Range('E47').Value = '  +1.80%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.553'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.63%  '
$ws.Range('E49').Value = '  +0.77%  '
$ws.Range('E50').Value = '  +0.57%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '10.81'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.41%  '
